$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Model Accuracy (-0.45, 0.45, 0.45)"
#   - new columns C:G (Market threshold, Market min, Market max, Recall,
#     Precision) with headers matching the style of the existing B1 header
#   - updated Accuracy (%) values in column B
# ---------------------------------------------------------------------------
$wsAcc = $wb.Worksheets.Item("Model Accuracy (-0.45, 0.45, 0.45)")

$wsAcc.Range("C1").Value = "Market threshold"
$wsAcc.Range("D1").Value = "Market min"
$wsAcc.Range("E1").Value = "Market max"
$wsAcc.Range("F1").Value = "Recall"
$wsAcc.Range("G1").Value = "Precision"

# Copy B1's header formatting (bold, centered, bordered) onto the new headers.
$wsAcc.Range("B1").Copy()
$wsAcc.Range("C1:G1").PasteSpecial(-4122)

# Row 2 - TOTALENERGIES SE
$wsAcc.Range("B2").Value = 59.65770171149144
$wsAcc.Range("C2").Value = 0.05450546436368681
$wsAcc.Range("D2").Value = -15.55441
$wsAcc.Range("E2").Value = 15.06418
$wsAcc.Range("F2").Value = 0
$wsAcc.Range("G2").Value = 0

# Row 3 - FMC CORP
$wsAcc.Range("B3").Value = 34.16870415647922
$wsAcc.Range("C3").Value = 0.009583939973006913
$wsAcc.Range("D3").Value = -19.35264
$wsAcc.Range("E3").Value = 13.70093
$wsAcc.Range("F3").Value = 5.630026809651475
$wsAcc.Range("G3").Value = 22.34042553191489

# Row 4 - BP PLC
$wsAcc.Range("B4").Value = 90.52567237163814
$wsAcc.Range("C4").Value = 0.04158117063764853
$wsAcc.Range("D4").Value = -18.75314
$wsAcc.Range("E4").Value = 23.33066
$wsAcc.Range("F4").Value = 0
$wsAcc.Range("G4").Value = 0

# Row 5 - STORA ENSO
$wsAcc.Range("B5").Value = 79.21760391198043
$wsAcc.Range("C5").Value = 0.02983403801513819
$wsAcc.Range("D5").Value = -12.78028
$wsAcc.Range("E5").Value = 12.42348
$wsAcc.Range("F5").Value = 0.9090909090909091
$wsAcc.Range("G5").Value = 5.555555555555555

# Row 6 - BHP GROUP
$wsAcc.Range("B6").Value = 92.05378973105135
$wsAcc.Range("C6").Value = 0.08368817696170747
$wsAcc.Range("D6").Value = -16.47904
$wsAcc.Range("E6").Value = 14.94325
$wsAcc.Range("F6").Value = 0
$wsAcc.Range("G6").Value = 0

# ---------------------------------------------------------------------------
# Sheet 2: "Confusion Matrix TOTALENERGIES SE (-0.45, 0.45, 0.45)"
# ---------------------------------------------------------------------------
$wsCm1 = $wb.Worksheets.Item("Confusion Matrix TOTALENERGIES SE (-0.45, 0.45, 0.45)")
$wsCm1.Range("B3").Value = 9
$wsCm1.Range("C3").Value = 973
$wsCm1.Range("D3").Value = 7

# ---------------------------------------------------------------------------
# Sheet 3: "Confusion Matrix FMC CORP (-0.45, 0.45, 0.45)"
# ---------------------------------------------------------------------------
$wsCm2 = $wb.Worksheets.Item("Confusion Matrix FMC CORP (-0.45, 0.45, 0.45)")
$wsCm2.Range("B2").Value = 21
$wsCm2.Range("C2").Value = 52
$wsCm2.Range("D2").Value = 21

$wsCm2.Range("B3").Value = 265
$wsCm2.Range("C3").Value = 464
$wsCm2.Range("D3").Value = 261

$wsCm2.Range("B4").Value = 87
$wsCm2.Range("C4").Value = 136
$wsCm2.Range("D4").Value = 74

# ---------------------------------------------------------------------------
# Sheet 4: "Confusion Matrix BP PLC (-0.45, 0.45, 0.45)"
# ---------------------------------------------------------------------------
$wsCm3 = $wb.Worksheets.Item("Confusion Matrix BP PLC (-0.45, 0.45, 0.45)")
$wsCm3.Range("B3").Value = 38
$wsCm3.Range("C3").Value = 1480
$wsCm3.Range("D3").Value = 41

$wsCm3.Range("B4").Value = 2
$wsCm3.Range("C4").Value = 36

# ---------------------------------------------------------------------------
# Sheet 5: "Confusion Matrix STORA ENSO (-0.45, 0.45, 0.45)"
# ---------------------------------------------------------------------------
$wsCm4 = $wb.Worksheets.Item("Confusion Matrix STORA ENSO (-0.45, 0.45, 0.45)")
$wsCm4.Range("B2").Value = 1
$wsCm4.Range("C2").Value = 16

$wsCm4.Range("B3").Value = 104
$wsCm4.Range("C3").Value = 1293
$wsCm4.Range("D3").Value = 104

$wsCm4.Range("B4").Value = 5
$wsCm4.Range("C4").Value = 49
$wsCm4.Range("D4").Value = 2

# ---------------------------------------------------------------------------
# Sheet 6: "Confusion Matrix BHP GROUP (-0.45, 0.45, 0.45)"
# ---------------------------------------------------------------------------
$wsCm5 = $wb.Worksheets.Item("Confusion Matrix BHP GROUP (-0.45, 0.45, 0.45)")
$wsCm5.Range("B2").Value = 0
$wsCm5.Range("C2").Value = 39

$wsCm5.Range("B3").Value = 4
$wsCm5.Range("C3").Value = 1506
$wsCm5.Range("D3").Value = 3
